$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Step 1: Replace the closing "Read our review ..." (italic) paragraph's
# text with the new DALLE image prompt text. Do this first, while the
# phrase is still unique in the document (before we add the new "Meta
# description" paragraph near the top, which also contains a copy of
# this same sentence).
# ---------------------------------------------------------------------
$oldClosing = "Read our review of Cirque du Soleil Kooza slot game, and play for free! Engaging gameplay mechanics with more win combos, colorful visual design and multiple special symbols and features."
$newClosing = "Prompt: DALLE, please create a feature image for Cirque du Soleil Kooza that captures the whimsical and colorful nature of the game and its circus theme. The image should be in a cartoon style and feature a happy Maya warrior with glasses. Make sure it is eye-catching and reflects the excitement and fun of this slot game."
$d.Content.Find.Execute($oldClosing, $false, $false, $false, $false, $false, $true, 1, $false, $newClosing, 2)

# ---------------------------------------------------------------------
# Step 2: Remove the duplicate bold "Play Cirque du Soleil Kooza for
# Free - Slot Game Review" paragraph that sits right before the
# (now updated) closing paragraph. Walk backwards and skip the very
# first paragraph (the real Heading1 title) so only the later, stray
# duplicate copy is removed.
# ---------------------------------------------------------------------
for ($i = $d.Paragraphs.Count; $i -ge 2; $i--) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -like "*Play Cirque du Soleil Kooza for Free - Slot Game Review*") {
        $p.Range.Delete()
    }
}

# ---------------------------------------------------------------------
# Step 3: Insert a new "Meta description" paragraph right after the
# opening Heading1 title paragraph.
# ---------------------------------------------------------------------
$title = $d.Paragraphs.First
$title.Range.InsertParagraphAfter()
$metaPara = $d.Paragraphs(2)
$metaXml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Meta description</w:t></w:r><w:r><w:t>: Read our review of Cirque du Soleil Kooza slot game, and play for free! Engaging gameplay mechanics with more win combos, colorful visual design and multiple special symbols and features.</w:t></w:r></w:p>"
$metaPara.Range.InsertXML($metaXml)

Write-Host "Edits applied."
